# Edit the workbook per the target diff:
# 1. Update the header texts in A1/B1/C1 (punctuation-only change: "16.5.1.1a." -> "16.5.1.1a ")
# 2. Add a new column I with 2020 data, matching the formatting of column H
# 3. Update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text updates (A1 / B1 / C1) ---
$ws.Range("A1").Value = '16.5.1.1a "Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк" индекси'
$ws.Range("B1").Value = '16.5.1.1a Индекс "Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления'''''
$ws.Range("C1").Value = '16.5.1.1a Index "Personal views about the level of corruption in executive government authorities and local government'''''

# --- 2. New column I (year 2020) ---

# Header cell I4 - same formatting as H4
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020

# I5 - same formatting as H5, with 0.0 number format
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 12.3
$ws.Range("I5").NumberFormat = "0.0"

# I6:I13 - same formatting as H6:H13 (borderless data rows), with 0.0 number format
$ws.Range("H6").Copy()
$ws.Range("I6:I13").PasteSpecial(-4122)
$ws.Range("I6:I13").NumberFormat = "0.0"
$ws.Range("I6").Value = 40.3
$ws.Range("I7").Value = 36.2
$ws.Range("I8").Value = 44.3
$ws.Range("I9").Value = 36
$ws.Range("I10").Value = 2.7
$ws.Range("I11").Value = 32.9
$ws.Range("I12").Value = 11.3
$ws.Range("I13").Value = -18.2

# I14 - same formatting as H14 (bottom-bordered), with 0.0 number format
$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 33
$ws.Range("I14").NumberFormat = "0.0"

# --- 3. Update selection ---
$ws.Range("F16").Select()
